$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# Finds the paragraph whose text starts with $startsWith; falls back to the
# 1-based $fallbackIndex if no (unique) match is found, so the script keeps
# working even if the document's paragraph ordering ever shifts slightly.
function Find-ParaByStart($startsWith, $fallbackIndex) {
    $match = $null
    $count = 0
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $txt = $d.Paragraphs($i).Range.Text
        if ($txt.StartsWith($startsWith)) {
            $count = $count + 1
            if ($match -eq $null) { $match = $i }
        }
    }
    if ($count -eq 1) { return $match }
    return $fallbackIndex
}

function Set-ParaXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs($paraIndex)
    $xml = $pkgHeader + '<w:body><w:p>' + $innerXml + '</w:p></w:body>' + $pkgFooter
    [void]$p.Range.InsertXML($xml)
}

$boldRPr = '<w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$szRPr = '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

# Paragraph 2: main.py
$inner = '<w:r>' + $boldRPr + '<w:t>main.py</w:t></w:r>' `
       + '<w:r>' + $szRPr + '<w:t xml:space="preserve"> </w:t></w:r>' `
       + '<w:r><w:br/></w:r>' `
       + '<w:r><w:t>Creates the GUI (through tkinter) and controls every other script</w:t></w:r>' `
       + '<w:r><w:t xml:space="preserve"> apart from Window.py</w:t></w:r>' `
       + '<w:r><w:t>.</w:t></w:r>'
Set-ParaXml (Find-ParaByStart "main.py" 2) $inner

# Paragraph 3: MazeDatabase.py
$inner = '<w:r>' + $boldRPr + '<w:t>MazeDatabase.py</w:t></w:r>' `
       + '<w:r>' + $szRPr + '<w:t xml:space="preserve">  </w:t></w:r>' `
       + '<w:r><w:br/></w:r>' `
       + '<w:r><w:t>Script for interaction with the maze.db for storing users and completed levels, uses SQL.</w:t></w:r>'
Set-ParaXml (Find-ParaByStart "MazeDatabase.py" 3) $inner

# Paragraph 4: MazeGenerationNew.py
$inner = '<w:r>' + $boldRPr + '<w:t>MazeGenerationNew.py</w:t></w:r>' `
       + '<w:r>' + $szRPr + '<w:t xml:space="preserve"> </w:t></w:r>' `
       + '<w:r>' + $szRPr + '<w:br/></w:r>' `
       + '<w:r><w:t>Uses Kruskal’s algorithm paired with a generated weight array to create the maze.</w:t></w:r>'
Set-ParaXml (Find-ParaByStart "MazeGenerationNew.py" 4) $inner

# Paragraph 5: MazeRendererNew.py
$inner = '<w:r>' + $boldRPr + '<w:t>MazeRendererNew.py</w:t></w:r>' `
       + '<w:r><w:br/></w:r>' `
       + '<w:r><w:t>Used when actually playing a maze, mainly used to implement the Window.py script</w:t></w:r>' `
       + '<w:r><w:t xml:space="preserve"> to correctly draw the maze</w:t></w:r>' `
       + '<w:r><w:t xml:space="preserve"> but also handles collision detection and a check to see if the player has won.</w:t></w:r>'
Set-ParaXml (Find-ParaByStart "MazeRendererNew.py" 5) $inner

# Paragraph 6: Window.py
$inner = '<w:r>' + $boldRPr + '<w:t>Window.py</w:t></w:r>' `
       + '<w:r>' + $szRPr + '<w:t xml:space="preserve"> </w:t></w:r>' `
       + '<w:r><w:br/></w:r>' `
       + '<w:r><w:t xml:space="preserve">Uses tkinter, NumPy and </w:t></w:r>' `
       + '<w:r><w:t>p</w:t></w:r>' `
       + '<w:r><w:t xml:space="preserve">illow </w:t></w:r>' `
       + '<w:r><w:t>(</w:t></w:r>' `
       + '<w:r><w:t>a fork of PIL (Python Imaging Library)</w:t></w:r>' `
       + '<w:r><w:t xml:space="preserve">) to create layered editable images on </w:t></w:r>' `
       + '<w:r><w:t>screen.</w:t></w:r>' `
       + '<w:r><w:t xml:space="preserve"> Using tkinter it also handles input, unfortunately the window needs to be clicked on for input to be registered.</w:t></w:r>'
Set-ParaXml (Find-ParaByStart "Window.py" 6) $inner

Write-Host "Done applying code-explanation formatting edits."
